$wb = $excel.ActiveWorkbook

# --- "Cars" sheet (existing, sheet1): replace/extend data rows ---
$ws1 = $wb.Worksheets.Item("Cars")
$ws1.Range("A2").Value = "Ester"
$ws1.Range("B2").Value = 12
$ws1.Range("C2").Value = "Pantera"
$ws1.Range("A3").Value = "Astor"
$ws1.Range("B3").Value = 13
$ws1.Range("C3").Value = "Linclon"

# --- "Sheet1" (new): original tester data plus an extra row, placed after "Cars" ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Sheet1"
$ws2.Move($null, $wb.Worksheets.Item("Cars"))

$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "age"
$ws2.Range("C1").Value = "car"
$ws2.Range("A2").Value = "tester"
$ws2.Range("B2").Value = 30
$ws2.Range("C2").Value = "sss"
$ws2.Range("A3").Value = "tester2"
$ws2.Range("B3").Value = 42
$ws2.Range("C3").Value = "Esteem"

# --- "Car" (new): scraped web listing, placed after "Sheet1" ---
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "Car"
$ws3.Move($null, $wb.Worksheets.Item("Sheet1"))

$ws3 = $wb.Worksheets.Item("Car")
$ws3.Range("A1").Value = "title"
$ws3.Range("B1").Value = "year"
$ws3.Range("C1").Value = "price"
$ws3.Range("D1").Value = "location"
$ws3.Range("E1").Value = "phone"
$ws3.Range("A2").Value = "'0714745050"
$ws3.Range("B2").Value = "'2014"
$ws3.Range("C2").Value = "Rs 3,025,000"
$ws3.Range("D2").Value = "Avissawella,"
$ws3.Range("E2").Value = "'0714745050"
